# Apply updated Fitness (column C) values for the run_11.xlsx log sheet.
# Rows 2-24  (Generation 0-22):  7622 -> 7318
# Rows 25-34 (Generation 23-32): 7312 -> 7318
# Rows 35-95 (Generation 33-93): 7293 -> 7318
# Rows 96-126 (Generation 94-124): 7293 -> 7310

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 95; $r++) {
    $ws.Cells.Item($r, 3).Value = 7318
}

for ($r = 96; $r -le 126; $r++) {
    $ws.Cells.Item($r, 3).Value = 7310
}
